# "wrapping up test file audit"
#
# The optimization_parameters sheet had a stray/obsolete "Sheet" row
# (row 16: label "Sheet" with values 3 and 4) left over from an earlier
# edit. Remove it, which shifts the simulation_timepoints row (old row
# 17) up to become row 16. Finish by leaving the threshold_b sheet as
# the active/selected tab, matching the final state the workbook was
# saved in.

$wb = $excel.ActiveWorkbook

# --- optimization_parameters: delete the obsolete "Sheet" row ---
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate()
$wsOpt.Rows.Item(16).Delete()

# Re-select the row that shifted up into row 16 (was row 17,
# "simulation_timepoints") so the sheet's recorded selection reflects
# the post-delete layout.
$wsOpt.Rows.Item(16).EntireRow.Select()

# --- make threshold_b the active sheet/tab ---
$wsThresh = $wb.Worksheets.Item("threshold_b")
$wsThresh.Activate()
